# Edit Specific Aims Page
# The "Heading 1" paragraph style previously relied on the spacing
# inherited from the "Normal" base style (after = 180 twips / 9pt).
# This change pins the style's own "space after paragraph" value to 0,
# which writes an explicit <w:spacing w:after="0"/> into the style's
# <w:pPr>, overriding the inherited spacing.

$d = $word.ActiveDocument

$heading1 = $d.Styles.Item("Heading 1")
$heading1.ParagraphFormat.SpaceAfter = 0
